$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 text from "Hno" to "Sr"
$ws.Range("C2").Value = "Sr"

# Clear out row 3 (Pablo Romero) and row 4 entirely (contents + formatting)
$ws.Range("A3:C4").Clear()

# Move the active selection to C2
$ws.Range("C2").Select()
